$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for 2021-09-10 .. 2021-09-20, appended after the existing
# last row (374). Column A keeps the same date-formatted style as the
# rows above it, so copy that formatting down first.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)

$data = @(
    @(375, 44449, 0, 12, 129.2546316242999),
    @(376, 44450, 1, 10, 107.7121930202499),
    @(377, 44451, 0, 6, 64.62731581214993),
    @(378, 44452, 2, 6, 64.62731581214993),
    @(379, 44453, 0, 5, 53.85609651012494),
    @(380, 44454, 0, 5, 53.85609651012494),
    @(381, 44455, 4, 7, 75.39853511417492),
    @(382, 44456, 0, 7, 75.39853511417492),
    @(383, 44457, 0, 6, 64.62731581214993),
    @(384, 44458, 4, 10, 107.7121930202499),
    @(385, 44459, 0, 8, 86.16975441619991)
)

foreach ($r in $data) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}
